# Applies the "Struccter fixed, Docs Started" restructuring to Sheet1:
#  - Drops the old Farnell/Comet "Магазин" side labels for the opamp rows and
#    instead moves the opamp part numbers (with their hyperlinks) + descriptions
#    into columns E/F for rows 2 and 3.
#  - Re-purposes column B rows 2/3 for two new components (2n3904 Transistor,
#    SMD Resistors/Capacitors 1206).
#  - Keeps MCP6291-E/MS (with hyperlink) + its description on row 4.
#  - Replaces the old "Arduino Nano 33 IoT" row with a new "Raspberry Pi Pico"
#    row that has a long, wrapped description and a bigger row height.
#  - Keeps the LED row (row 6) as-is (hyperlink preserved).
#  - Adds a new "Shcottkey Diods" row 7.
#  - Normalizes alignment (drop horizontal centering, add vertical centering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove all existing hyperlinks (targets are re-created further down in
#    their new locations). Doing this before touching cell values means we
#    do not have to fight with stale hyperlink relationships pointing at
#    cells that are about to hold different content.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Clear out the old layout completely so we can rebuild it cleanly.
# ---------------------------------------------------------------------------
$ws.Range("A1:F9").ClearContents()

# ---------------------------------------------------------------------------
# 3. Header row (unchanged text/position, only alignment changes later).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Компоненти"
$ws.Range("C1").Value = "описание"
$ws.Range("D1").Value = "Магазин"

# ---------------------------------------------------------------------------
# 4. Row 2: new "2n3904 Transistor" component in B, old MCP6L91T-E/OT moved to
#    E/F with its hyperlink + description.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "2n3904 Transistor"
$ws.Range("E2").Value = "MCP6L91T-E/OT"
$ws.Range("F2").Value = "Operational Amplifier, Single, 1 Channels, 10 MHz, 7 V/µs, 2.4V to 6V, SOT-23, 5 Pins"

# ---------------------------------------------------------------------------
# 5. Row 3: new "SMD Resistors/Capacitors 1206" component in B, old
#    MCP6291T-E/OT moved to E/F with its hyperlink + description.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "SMD Resistors/Capacitors 1206"
$ws.Range("E3").Value = "MCP6291T-E/OT"
$ws.Range("F3").Value = "Operational Amplifier, Single, 1 Channels, 10 MHz, 7 V/µs, 2.4V to 5.5V, SOT-23, 5 Pins"

# ---------------------------------------------------------------------------
# 6. Row 4: MCP6291-E/MS keeps its hyperlink + description, moved up from row 3.
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "MCP6291-E/MS"
$ws.Range("C4").Value = "Operational Amplifier, 1 OA,3mV,1pA,2.4-6V,1mA,10MHz"

# ---------------------------------------------------------------------------
# 7. Row 5: replaces "Arduino Nano 33 IoT" with "Raspberry Pi Pico" plus a
#    long wrapped description (keeps the old hyperlink-blue font look even
#    though the hyperlink itself is gone).
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Raspberry Pi Pico"
$ws.Range("C5").Value = "Raspberry Pi Pico 2 is a low-cost, high-performance microcontroller board with flexible digital interfaces. Key features include:`nRP2350 `nDual Cortex-M33 up to 150MHz`n520KB of SRAM, 4MB flash memory`nUSB 1.1 with device and host support`nLow-power sleep and dormant modes`nDrag-and-drop programming using mass storage over USB`n26× GPIO pins, 3 ADC`n2× SPI, 2× I2C, 2× UART, 3× 12-bit 500ksps Analogue to Digital Converter (ADC), 24× controllable PWM channels`n2× Timer with 4 alarms, 1× AON Timer`nTemperature sensor"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("C5").WrapText = $true

# ---------------------------------------------------------------------------
# 8. Row 6: LED stays the same (hyperlink re-added below).
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "OSTB8BS4C2B LED"
$ws.Range("C6").Value = "5.0x5.0x1.5mm, "

# ---------------------------------------------------------------------------
# 9. Row 7: brand-new "Shcottkey Diods" component.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Shcottkey Diods"

# ---------------------------------------------------------------------------
# 10. Re-create the hyperlinks at their new locations.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E2"), "https://de.farnell.com/en-DE/microchip/mcp6l91t-e-ot/ic-op-amp-single-10mhz-5sot-23/dp/1715865")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://store.comet.bg/Catalogue/Product/16935/")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://de.farnell.com/en-DE/microchip/mcp6291t-e-ot/op-amp-10mhz-single-smd-sot235/dp/1439464")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://store.comet.bg/Catalogue/Product/5301166/")

# ---------------------------------------------------------------------------
# 11. Alignment clean-up: remove horizontal centering, add vertical centering
#     to every cell that is actually in use (matches the "Struccter fixed"
#     part of the commit).
# ---------------------------------------------------------------------------
$used = $ws.Range("B1:F9")
$used.VerticalAlignment = -4108   # xlVAlignCenter
$used.HorizontalAlignment = 1     # xlHAlignGeneral

# ---------------------------------------------------------------------------
# 12. Row height for the long Raspberry Pi Pico description & merged footer
#     row stays a plain merged/empty row (already merged in the template).
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 199.5

# ---------------------------------------------------------------------------
# 13. Column widths (best achievable approximation of the target widths; the
#     automation layer quantizes column widths to 1/6-character steps).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.166666666666666
$ws.Columns.Item(2).ColumnWidth = 32.5
$ws.Columns.Item(3).ColumnWidth = 71.16666666666667
$ws.Columns.Item(4).ColumnWidth = 24.166666666666668
$ws.Columns.Item(5).ColumnWidth = 14.333333333333332
$ws.Columns.Item(6).ColumnWidth = 30.5

# ---------------------------------------------------------------------------
# 14. Selection / active cell now sits on the merged footer row.
# ---------------------------------------------------------------------------
$ws.Range("B9:D9").Select()

Write-Output "edit complete"
